$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = 'Última actualización: 18:52:02'
$ws.Cells.Item(3, 1).Value = 'Total filas: 328'
$ws.Cells.Item(64, 1).Value = '08:27:16'
$ws.Cells.Item(64, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(64, 4).Value = 72
$ws.Cells.Item(65, 1).Value = '07:50:33'
$ws.Cells.Item(65, 3).Value = '15_ABASTO'
$ws.Cells.Item(65, 4).Value = 109
$ws.Cells.Item(86, 1).Value = '09:23:23'
$ws.Cells.Item(86, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(86, 4).Value = 96
$ws.Cells.Item(87, 1).Value = '10:50:41'
$ws.Cells.Item(87, 3).Value = '10_OLMOS'
$ws.Cells.Item(87, 4).Value = 9
$ws.Cells.Item(106, 1).Value = '11:52:01'
$ws.Cells.Item(106, 3).Value = '15X38_ABASTO'
$ws.Cells.Item(106, 4).Value = 0
$ws.Cells.Item(108, 1).Value = '11:47:17'
$ws.Cells.Item(108, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(108, 4).Value = 5
$ws.Cells.Item(118, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(119, 3).Value = '15_ABASTO'
$ws.Cells.Item(128, 1).Value = '11:47:17'
$ws.Cells.Item(128, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(128, 4).Value = 45
$ws.Cells.Item(129, 1).Value = '10:37:52'
$ws.Cells.Item(129, 3).Value = '14_ABASTO'
$ws.Cells.Item(129, 4).Value = 115
$ws.Cells.Item(137, 1).Value = '11:52:01'
$ws.Cells.Item(137, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(137, 4).Value = 45
$ws.Cells.Item(138, 1).Value = '11:47:17'
$ws.Cells.Item(138, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(138, 4).Value = 50
$ws.Cells.Item(167, 1).Value = '12:45:56'
$ws.Cells.Item(167, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(167, 4).Value = 76
$ws.Cells.Item(168, 1).Value = '12:11:52'
$ws.Cells.Item(168, 3).Value = '10_OLMOS'
$ws.Cells.Item(168, 4).Value = 110
$ws.Cells.Item(169, 1).Value = '12:33:21'
$ws.Cells.Item(169, 3).Value = '10_OLMOS'
$ws.Cells.Item(169, 4).Value = 89
$ws.Cells.Item(170, 1).Value = '13:14:29'
$ws.Cells.Item(170, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(170, 4).Value = 48
$ws.Cells.Item(209, 1).Value = '13:56:11'
$ws.Cells.Item(209, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(209, 4).Value = 117
$ws.Cells.Item(210, 3).Value = '15X38_ABASTO'
$ws.Cells.Item(211, 1).Value = '14:53:07'
$ws.Cells.Item(211, 3).Value = '10_OLMOS'
$ws.Cells.Item(211, 4).Value = 60
$ws.Cells.Item(301, 1).Value = '18:44:57'
$ws.Cells.Item(301, 3).Value = '14X44_ABASTO'
$ws.Cells.Item(301, 4).Value = 32
$ws.Cells.Item(302, 1).Value = '17:56:03'
$ws.Cells.Item(302, 3).Value = '17_ROMERO'
$ws.Cells.Item(302, 4).Value = 80
$ws.Cells.Item(303, 1).Value = '18:12:30'
$ws.Cells.Item(303, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(303, 4).Value = 64
$ws.Cells.Item(312, 1).Value = '18:52:02'
$ws.Cells.Item(312, 2).Value = '19:29'
$ws.Cells.Item(312, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(312, 4).Value = 37
$ws.Cells.Item(313, 1).Value = '17:48:33'
$ws.Cells.Item(313, 2).Value = '19:35'
$ws.Cells.Item(313, 4).Value = 107
$ws.Cells.Item(314, 1).Value = '17:56:03'
$ws.Cells.Item(314, 2).Value = '19:36'
$ws.Cells.Item(314, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(314, 4).Value = 100
$ws.Cells.Item(315, 1).Value = '18:44:57'
$ws.Cells.Item(315, 2).Value = '19:38'
$ws.Cells.Item(315, 4).Value = 54
$ws.Cells.Item(316, 1).Value = '17:48:33'
$ws.Cells.Item(316, 2).Value = '19:39'
$ws.Cells.Item(316, 3).Value = '15X38_ABASTO'
$ws.Cells.Item(316, 4).Value = 111
$ws.Cells.Item(317, 1).Value = '18:44:57'
$ws.Cells.Item(317, 2).Value = '19:51'
$ws.Cells.Item(317, 4).Value = 67
$ws.Cells.Item(318, 1).Value = '17:56:03'
$ws.Cells.Item(318, 3).Value = '81_EL PELIGRO'
$ws.Cells.Item(318, 4).Value = 116
$ws.Cells.Item(319, 1).Value = '18:44:57'
$ws.Cells.Item(319, 2).Value = '19:52'
$ws.Cells.Item(319, 4).Value = 68
$ws.Cells.Item(320, 1).Value = '17:56:03'
$ws.Cells.Item(320, 3).Value = '225_GOMEZ'
$ws.Cells.Item(320, 4).Value = 117
$ws.Cells.Item(321, 2).Value = '19:53'
$ws.Cells.Item(321, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(321, 4).Value = 69
$ws.Cells.Item(322, 2).Value = '20:06'
$ws.Cells.Item(322, 3).Value = '215C_EL PATO'
$ws.Cells.Item(322, 4).Value = 82
$ws.Cells.Item(323, 1).Value = '18:52:02'
$ws.Cells.Item(323, 2).Value = '20:07'
$ws.Cells.Item(323, 4).Value = 75
$ws.Cells.Item(324, 1).Value = '18:52:02'
$ws.Cells.Item(324, 2).Value = '20:08'
$ws.Cells.Item(324, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(324, 4).Value = 76
$ws.Cells.Item(325, 2).Value = '20:09'
$ws.Cells.Item(325, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(325, 4).Value = 85
$ws.Cells.Item(326, 1).Value = '18:44:57'
$ws.Cells.Item(326, 2).Value = '20:12'
$ws.Cells.Item(326, 3).Value = '14_ABASTO'
$ws.Cells.Item(326, 4).Value = 88
$ws.Cells.Item(327, 1).Value = '18:31:25'
$ws.Cells.Item(327, 2).Value = '20:12'
$ws.Cells.Item(327, 3).Value = '215C_EL PATO'
$ws.Cells.Item(327, 4).Value = 101
$ws.Cells.Item(328, 1).Value = '18:52:02'
$ws.Cells.Item(328, 2).Value = '20:13'
$ws.Cells.Item(328, 3).Value = '14_ABASTO'
$ws.Cells.Item(328, 4).Value = 81
$ws.Cells.Item(328, 5).Value = 'LP1912'
$ws.Cells.Item(329, 1).Value = '18:44:57'
$ws.Cells.Item(329, 2).Value = '20:21'
$ws.Cells.Item(329, 3).Value = '15_ABASTO'
$ws.Cells.Item(329, 4).Value = 97
$ws.Cells.Item(329, 5).Value = 'LP1912'
$ws.Cells.Item(330, 1).Value = '18:31:25'
$ws.Cells.Item(330, 2).Value = '20:22'
$ws.Cells.Item(330, 3).Value = '15_ABASTO'
$ws.Cells.Item(330, 4).Value = 111
$ws.Cells.Item(330, 5).Value = 'LP1912'
$ws.Cells.Item(331, 1).Value = '18:44:57'
$ws.Cells.Item(331, 2).Value = '20:30'
$ws.Cells.Item(331, 3).Value = '10_OLMOS'
$ws.Cells.Item(331, 4).Value = 106
$ws.Cells.Item(331, 5).Value = 'LP1912'
$ws.Cells.Item(332, 1).Value = '18:52:02'
$ws.Cells.Item(332, 2).Value = '20:31'
$ws.Cells.Item(332, 3).Value = '10_OLMOS'
$ws.Cells.Item(332, 4).Value = 99
$ws.Cells.Item(332, 5).Value = 'LP1912'
$ws.Cells.Item(333, 1).Value = '18:52:02'
$ws.Cells.Item(333, 2).Value = '20:48'
$ws.Cells.Item(333, 3).Value = '215B_EL PATO'
$ws.Cells.Item(333, 4).Value = 116
$ws.Cells.Item(333, 5).Value = 'LP1912'

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = 'Última actualización: 18:52:02'
$ws.Cells.Item(3, 1).Value = 'Total filas: 53'
$ws.Cells.Item(56, 1).Value = '18:52:02'
$ws.Cells.Item(56, 2).Value = '20:07'
$ws.Cells.Item(56, 4).Value = 75
$ws.Cells.Item(57, 1).Value = '18:31:25'
$ws.Cells.Item(57, 2).Value = '20:12'
$ws.Cells.Item(57, 3).Value = '215C_EL PATO'
$ws.Cells.Item(57, 4).Value = 101
$ws.Cells.Item(57, 5).Value = 'LP1912'
$ws.Cells.Item(58, 1).Value = '18:52:02'
$ws.Cells.Item(58, 2).Value = '20:48'
$ws.Cells.Item(58, 3).Value = '215B_EL PATO'
$ws.Cells.Item(58, 4).Value = 116
$ws.Cells.Item(58, 5).Value = 'LP1912'

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = 'Última actualización: 18:52:02'
